$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Un-minimize the workbook window
$wb.Windows.Item(1).WindowState = [Microsoft.Office.Interop.Excel.XlWindowState]::xlNormal

# Add the two new shared strings via a new row 73
$ws.Cells.Item(73, 1).Value = "BASIC_TEXT_OFFICERS_REQUIRED"
$ws.Cells.Item(73, 2).Value = "Officers required"
$ws.Cells.Item(73, 3).Value = "XXXX"
$ws.Cells.Item(73, 4).Value = "XXXX"
$ws.Cells.Item(73, 5).Value = "XXXX"

# Update the view: scroll position and selection
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E73").Select()
